$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: it becomes a "continuation" row (like rows 3 and 6) of the
# previous entry, so it gets the bordered/no-top continuation style and an
# (empty) A8 cell gains that same style too ---
$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A8:E8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Add the two new data rows (9 and 10) ---
# Copy the formatting of row 7 (a normal/"main" entry row) onto the new rows
# first, then fill in the values. The order in which string values are
# assigned below matches the order new shared strings were appended in the
# original authoring session (filename, EN x2, RU x2, corrupted x2, filename).
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A9:E9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A10:E10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A9").Value = "SCRIPT/P01P04A/um1302.ssb"
$ws.Range("C9").Value = " I am back from the\nlonely expedition!"
$ws.Range("C10").Value = " It was such a dangerous\nadventure, but I feel like it let me gain so\nmuch strength!"
$ws.Range("D9").Value = " Я вернулся из одиночной\nэкспедиции!"
$ws.Range("D10").Value = " Это было очень опасное\nприключение, но я чувствую, что я стал\nгораздо сильнее!"
$ws.Range("E9").Value = " Ÿ âåñîôìòÿ éè ïäéîïœîïê\nüëòðåäéøéé!"
$ws.Range("E10").Value = " Üóï áúìï ïœåîû ïðàòîïå\nðñéëìýœåîéå, îï ÿ œôâòóâôý, œóï ÿ òóàì\nãïñàèäï òéìûîåå!"
$ws.Range("A10").Value = "SCRIPT/P01P04A/um1408.ssb"

$ws.Range("B9").Value = 267
$ws.Range("B10").Value = 270

# --- Row heights for the two new rows (wrapped multi-line text) ---
$ws.Rows.Item(9).RowHeight = 43.2
$ws.Rows.Item(10).RowHeight = 43.2

# --- View state: scrolled down a bit further, selection moved to C10 ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C10").Select() | Out-Null
